# Generate Report for Archive
#
# 1) The shared status string "Ready for handoff" becomes "In Translation"
#    everywhere it is used (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) The "Status" column is narrowed on all three sheets (Overview columns
#    E & F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns (was ~17.22 OOXML chars wide, now ~13.41).
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
